$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (154-167) to append: Fecha, Microciclo_Num, Tipo_Microciclo, Fase, Tipo_Dia, Intensidad, Partido
$rows = @(
    @{ r = 154; fecha = 45950; micro = 23; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = 2;    partido = $null },
    @{ r = 155; fecha = 45951; micro = 23; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = -1;   partido = $null },
    @{ r = 156; fecha = 45952; micro = 23; tipo = "Competencia"; fase = "Competencia"; dia = "PARTIDO";   intensidad = $null; partido = "Queretaro" },
    @{ r = 157; fecha = 45953; micro = 23; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = 1;    partido = $null },
    @{ r = 158; fecha = 45954; micro = 23; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = -1;   partido = $null },
    @{ r = 159; fecha = 45955; micro = 23; tipo = "Competencia"; fase = "Competencia"; dia = "PARTIDO";   intensidad = $null; partido = "Atlas" },
    @{ r = 160; fecha = 45956; micro = 23; tipo = "Competencia"; fase = "Competencia"; dia = "DESCANSO";  intensidad = $null; partido = $null },
    @{ r = 161; fecha = 45957; micro = 24; tipo = "Competencia"; fase = "Competencia"; dia = "DESCANSO";  intensidad = $null; partido = $null },
    @{ r = 162; fecha = 45958; micro = 24; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = 1;    partido = $null },
    @{ r = 163; fecha = 45959; micro = 24; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = 2;    partido = $null },
    @{ r = 164; fecha = 45960; micro = 24; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = 3;    partido = $null },
    @{ r = 165; fecha = 45961; micro = 24; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = -2;   partido = $null },
    @{ r = 166; fecha = 45962; micro = 24; tipo = "Competencia"; fase = "Competencia"; dia = "ENTRENO";   intensidad = -1;   partido = $null },
    @{ r = 167; fecha = 45963; micro = 24; tipo = "Competencia"; fase = "Competencia"; dia = "PARTIDO";   intensidad = $null; partido = "Pachuca" }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value = $row.fecha
    $ws.Cells.Item($row.r, 2).Value = $row.micro
    $ws.Cells.Item($row.r, 3).Value = $row.tipo
    $ws.Cells.Item($row.r, 4).Value = $row.fase
    $ws.Cells.Item($row.r, 5).Value = $row.dia
    if ($null -ne $row.intensidad) {
        $ws.Cells.Item($row.r, 6).Value = $row.intensidad
    }
    if ($null -ne $row.partido) {
        $ws.Cells.Item($row.r, 7).Value = $row.partido
    }
}

# Update the view state to match the final saved workbook
$ws.Application.ActiveWindow.ScrollRow = 97
$ws.Range("H165").Select()
